$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("Factory"), shifting all
# existing data rows down by one (row 2 -> 3, 3 -> 4, ... 6 -> 7).
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new "salary" income entry.
$ws.Range("A2").Value = "salary"
$ws.Range("C2").Value = 45846.59566582176

# Copy the number formatting (date style) from the row below so the new
# date cell matches the existing date column formatting/style.
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
